$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 19.571477
$ws.Range("H2").Value = 58.714431
$ws.Range("I2").Value = 0.07340284417718584
$ws.Range("J2").Value = 0.07340284417718584
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.7941993333333333
$ws.Range("N2").Value = 2.382598
$ws.Range("O2").Value = 0.475568178302488
$ws.Range("P2").Value = 0.475568178302488
$ws.Range("Q2").Value = 15.54365398574866
$ws.Range("R2").Value = 139.892885871738
$ws.Range("S2").Value = 0.03490805688756567
$ws.Range("T2").Value = 0.03490805688756566

# Row 3
$ws.Range("G3").Value = 19.571477
$ws.Range("H3").Value = 58.714431
$ws.Range("I3").Value = 0.07340284417718584
$ws.Range("J3").Value = 0.07340284417718584
$ws.Range("O3").Value = 0.4790747630290841
$ws.Range("P3").Value = 0.479074763029084
$ws.Range("Q3").Value = 15.65826455506066
$ws.Range("R3").Value = 140.924380995546
$ws.Range("S3").Value = 0.03516545017984609
$ws.Range("T3").Value = 0.03516545017984608

# Row 4
$ws.Range("G4").Value = 19.571477
$ws.Range("H4").Value = 58.714431
$ws.Range("I4").Value = 0.07340284417718584
$ws.Range("J4").Value = 0.07340284417718584
$ws.Range("O4").Value = 0.04535705866842794
$ws.Range("P4").Value = 0.04535705866842794
$ws.Range("Q4").Value = 1.482467620667667
$ws.Range("R4").Value = 13.342208586009
$ws.Range("S4").Value = 0.003329337109774093
$ws.Range("T4").Value = 0.003329337109774093

# Row 5
$ws.Range("I5").Value = 0.9221772503952724
$ws.Range("J5").Value = 0.9221772503952725
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.7941993333333333
$ws.Range("N5").Value = 2.382598
$ws.Range("O5").Value = 0.475568178302488
$ws.Range("P5").Value = 0.475568178302488
$ws.Range("Q5").Value = 195.2785924626111
$ws.Range("R5").Value = 1757.5073321635
$ws.Range("S5").Value = 0.4385581550424771
$ws.Range("T5").Value = 0.4385581550424771

# Row 6
$ws.Range("I6").Value = 0.9221772503952724
$ws.Range("J6").Value = 0.9221772503952725
$ws.Range("O6").Value = 0.4790747630290841
$ws.Range("P6").Value = 0.479074763029084
$ws.Range("S6").Value = 0.4417918477039275
$ws.Range("T6").Value = 0.4417918477039274

# Row 7
$ws.Range("I7").Value = 0.9221772503952724
$ws.Range("J7").Value = 0.9221772503952725
$ws.Range("O7").Value = 0.04535705866842794
$ws.Range("P7").Value = 0.04535705866842794
$ws.Range("S7").Value = 0.04182724764886794
$ws.Range("T7").Value = 0.04182724764886794

# Row 8
$ws.Range("I8").Value = 0.004419905427541656
$ws.Range("J8").Value = 0.004419905427541656
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.7941993333333333
$ws.Range("N8").Value = 2.382598
$ws.Range("O8").Value = 0.475568178302488
$ws.Range("P8").Value = 0.475568178302488
$ws.Range("Q8").Value = 0.9359512071439998
$ws.Range("R8").Value = 8.423560864295998
$ws.Range("S8").Value = 0.002101966372445265
$ws.Range("T8").Value = 0.002101966372445265

# Row 9
$ws.Range("I9").Value = 0.004419905427541656
$ws.Range("J9").Value = 0.004419905427541656
$ws.Range("O9").Value = 0.4790747630290841
$ws.Range("P9").Value = 0.479074763029084
$ws.Range("Q9").Value = 0.9428524094479999
$ws.Range("R9").Value = 8.485671685031999
$ws.Range("S9").Value = 0.002117465145310481
$ws.Range("T9").Value = 0.002117465145310481

# Row 10
$ws.Range("I10").Value = 0.004419905427541656
$ws.Range("J10").Value = 0.004419905427541656
$ws.Range("O10").Value = 0.04535705866842794
$ws.Range("P10").Value = 0.04535705866842794
$ws.Range("Q10").Value = 0.08926584189199999
$ws.Range("R10").Value = 0.8033925770279999
$ws.Range("S10").Value = 0.00020047390978591
$ws.Range("T10").Value = 0.00020047390978591
